$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 115
$ws.Range("J29").Value = 115
$ws.Range("L29").Value = 345
$ws.Range("N29").Value = -907
$ws.Range("H116").Value = 9397.333000000001
$ws.Range("I116").Value = 6260.6665
$ws.Range("J116").Value = 10965.667
$ws.Range("K116").Value = 6260.6665
$ws.Range("L116").Value = 10965.667
$ws.Range("M116").Value = -2818.6665
$ws.Range("N116").Value = -17849.667
$ws.Range("H127").Value = 1999.5
$ws.Range("I127").Value = 1999.5
$ws.Range("K127").Value = 5998.5
$ws.Range("M127").Value = -1038.5
$ws.Range("H132").Value = 1862.125
$ws.Range("I132").Value = 1862.125
$ws.Range("K132").Value = 5586.375
$ws.Range("M132").Value = -3056.375
$ws.Range("H137").Value = 3728.8909
$ws.Range("J137").Value = 8684.588
$ws.Range("L137").Value = 26053.764
$ws.Range("N137").Value = -31153.764
$ws.Range("H138").Value = 3701.2188
$ws.Range("J138").Value = 3731.25
$ws.Range("L138").Value = 11193.75
$ws.Range("N138").Value = -21473.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 8971.143
$ws.Range("J21").Value = 9574.75
$ws.Range("L21").Value = 9574.75
$ws.Range("N21").Value = -10322.75
$ws.Range("H45").Value = 13245
$ws.Range("I45").Value = 16874.375
$ws.Range("J45").Value = 3566.6667
$ws.Range("K45").Value = 16874.375
$ws.Range("L45").Value = 3566.6667
$ws.Range("M45").Value = -16497.375
$ws.Range("N45").Value = -4320.6667
$ws.Range("H61").Value = 4002.6223
$ws.Range("I61").Value = 3011.543
$ws.Range("K61").Value = 3011.543
$ws.Range("M61").Value = -2799.543
$ws.Range("H74").Value = 3736.9656
$ws.Range("I74").Value = 3296.6
$ws.Range("K74").Value = 3296.6
$ws.Range("M74").Value = -2422.6
$ws.Range("H77").Value = 3736.9656
$ws.Range("I77").Value = 3296.6
$ws.Range("K77").Value = 16483
$ws.Range("M77").Value = -12115
$ws.Range("H110").Value = 1671.1904
$ws.Range("I110").Value = 1227.5555
$ws.Range("K110").Value = 1227.5555
$ws.Range("M110").Value = 817.4445000000001
$ws.Range("H122").Value = 1786.762
$ws.Range("I122").Value = 1704
$ws.Range("J122").Value = 2283.3333
$ws.Range("K122").Value = 5112
$ws.Range("L122").Value = 6849.999899999999
$ws.Range("M122").Value = -2662
$ws.Range("N122").Value = -11749.9999
$ws.Range("H132").Value = 2878.3572
$ws.Range("I132").Value = 2902.6667
$ws.Range("J132").Value = 2222
$ws.Range("K132").Value = 8708.000100000001
$ws.Range("L132").Value = 6666
$ws.Range("M132").Value = -6178.000100000001
$ws.Range("N132").Value = -11726
$ws.Range("H136").Value = 4002.6223
$ws.Range("I136").Value = 3011.543
$ws.Range("K136").Value = 9034.629000000001
$ws.Range("M136").Value = -6484.629000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 272298.8
$ws.Range("I86").Value = 372161.88
$ws.Range("J86").Value = 2668.5
$ws.Range("K86").Value = 372161.88
$ws.Range("L86").Value = 2668.5
$ws.Range("M86").Value = -371038.88
$ws.Range("N86").Value = -4914.5
$ws.Range("H88").Value = 15371
$ws.Range("J88").Value = 15371
$ws.Range("L88").Value = 15371
$ws.Range("N88").Value = -16183
$ws.Range("H89").Value = 272298.8
$ws.Range("I89").Value = 372161.88
$ws.Range("J89").Value = 2668.5
$ws.Range("K89").Value = 1860809.4
$ws.Range("L89").Value = 13342.5
$ws.Range("M89").Value = -1855193.4
$ws.Range("N89").Value = -24574.5
$ws.Range("H91").Value = 15371
$ws.Range("J91").Value = 15371
$ws.Range("L91").Value = 15371
$ws.Range("N91").Value = -18179
$ws.Range("H94").Value = 857.63635
$ws.Range("I94").Value = 839.375
$ws.Range("J94").Value = 906.3333
$ws.Range("K94").Value = 839.375
$ws.Range("L94").Value = 906.3333
$ws.Range("M94").Value = -388.375
$ws.Range("N94").Value = -1808.3333
$ws.Range("H105").Value = 2439.8696
$ws.Range("I105").Value = 2536.4736
$ws.Range("J105").Value = 2371.889
$ws.Range("K105").Value = 2536.4736
$ws.Range("L105").Value = 2371.889
$ws.Range("M105").Value = -789.4735999999998
$ws.Range("N105").Value = -5865.889
$ws.Range("H134").Value = 16707.25
$ws.Range("I134").Value = 5955.7144
$ws.Range("K134").Value = 17867.1432
$ws.Range("M134").Value = -15332.1432
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 13339965
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 13339965
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = 13339965
$ws.Range("N4").Value = -13340189
$ws.Range("H105").Value = 2564.6667
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 4328.5415
$ws.Range("I132").Value = 3812.9546
$ws.Range("K132").Value = 11438.8638
$ws.Range("M132").Value = -8908.863799999999
$ws.Range("H134").Value = 6289.2
$ws.Range("I134").Value = 6891.5
$ws.Range("J134").Value = 3880
$ws.Range("K134").Value = 20674.5
$ws.Range("L134").Value = 11640
$ws.Range("M134").Value = -18139.5
$ws.Range("N134").Value = -16710
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3240.2
$ws.Range("I129").Value = 1507.8
$ws.Range("J129").Value = 4972.6
$ws.Range("K129").Value = 4523.4
$ws.Range("L129").Value = 14917.8
$ws.Range("M129").Value = 476.6000000000004
$ws.Range("N129").Value = -24917.8
$ws.Range("H131").Value = 14515.77
$ws.Range("J131").Value = 1848.4395
$ws.Range("L131").Value = 5545.318499999999
$ws.Range("N131").Value = -15625.3185
$ws.Range("H137").Value = 4994
$ws.Range("J137").Value = 5374.75
$ws.Range("L137").Value = 16124.25
$ws.Range("N137").Value = -26324.25
$ws.Range("H140").Value = 2835.8572
$ws.Range("I140").Value = 2391.8333
$ws.Range("K140").Value = 7175.499899999999
$ws.Range("M140").Value = -1995.499899999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 30992.666
$ws.Range("J98").Value = 30992.666
$ws.Range("L98").Value = 30992.666
$ws.Range("N98").Value = -36982.666
$ws.Range("H111").Value = 29000
$ws.Range("J111").Value = 29000
$ws.Range("L111").Value = 29000
$ws.Range("N111").Value = -35134
$ws.Range("H132").Value = 9051.799999999999
$ws.Range("J132").Value = 2042.4
$ws.Range("L132").Value = 6127.200000000001
$ws.Range("N132").Value = -11187.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("N14").Value = 0
$ws.Range("H38").Value = 25854.334
$ws.Range("I38").Value = 19230
$ws.Range("J38").Value = 29166.5
$ws.Range("K38").Value = 19230
$ws.Range("L38").Value = 29166.5
$ws.Range("M38").Value = -18820
$ws.Range("N38").Value = -29986.5
$ws.Range("H100").Value = 5798.5
$ws.Range("I100").Value = 4647.75
$ws.Range("K100").Value = 4647.75
$ws.Range("M100").Value = -4106.75
$ws.Range("H103").Value = 39998
$ws.Range("J103").Value = 39996
$ws.Range("L103").Value = 39996
$ws.Range("N103").Value = -42340
$ws.Range("H110").Value = 70883.336
$ws.Range("J110").Value = 70883.336
$ws.Range("L110").Value = 70883.336
$ws.Range("N110").Value = -79063.336
$ws.Range("H122").Value = 3168.3
$ws.Range("I122").Value = 2858.2856
$ws.Range("K122").Value = 8574.856800000001
$ws.Range("M122").Value = -6124.856800000001
$ws.Range("H132").Value = 2894.7693
$ws.Range("I132").Value = 2526.359
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7579.076999999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5049.076999999999
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 2372.2407
$ws.Range("I136").Value = 2072.02
$ws.Range("K136").Value = 6216.059999999999
$ws.Range("M136").Value = -3666.059999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 220.2
$ws.Range("J23").Value = 900
$ws.Range("L23").Value = 900
$ws.Range("N23").Value = -1358
$ws.Range("H37").Value = 46666
$ws.Range("I37").Value = 46666
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 46666
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -46463
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = 0
$ws.Range("H49").Value = 50000
$ws.Range("I49").Value = 50000
$ws.Range("K49").Value = 50000
$ws.Range("M49").Value = -49770
$ws.Range("H122").Value = 4731.75
$ws.Range("I122").Value = 2741.7144
$ws.Range("K122").Value = 8225.143199999999
$ws.Range("M122").Value = -5775.143199999999
$ws.Range("H136").Value = 1018.4
$ws.Range("I136").Value = 1051.9131
$ws.Range("K136").Value = 3155.7393
$ws.Range("M136").Value = -605.7393000000002
